$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gained two new sighting rows (dates 2023-05-06 and 2023-05-07),
# which got sorted into the existing, date-ordered table right before the last
# two rows. Concretely:
#   old row 49 (2023-05-07, Saint-Jean-sur-le-Richelieu)  -> becomes new row 50
#   old row 50 (2023-05-08, Vaudreuil)                     -> becomes new row 52
#   new row 49 (2023-05-06, Montreal)                      is inserted
#   new row 51 (2023-05-07, Montcalm)                      is inserted
#
# Push the existing data down first (bottom row first so we don't clobber
# data we still need to copy), reusing Range.Copy(Destination) so both the
# values AND the existing cell formatting/styles are carried over exactly.
$ws.Range("A50:I50").Copy($ws.Range("A52:I52"))
$ws.Range("A49:I49").Copy($ws.Range("A50:I50"))

# Row 51 is a brand-new row whose "Zone climatique" is "B" (the orange-shaded
# style used elsewhere in the sheet, e.g. row 29, which also shares the same
# "BUAM" contact styling as the other rows here), so clone formatting from
# that kind of row before writing its values.
$ws.Range("A29:I29").Copy($ws.Range("A51:I51"))

# Row 49 (new sighting, 2023-05-06 / Montreal)
$ws.Cells.Item(49, 1).Value = "5/6/2023"
$ws.Cells.Item(49, 2).Value = "BUAM"
$ws.Cells.Item(49, 3).Value = "N/A"
$ws.Cells.Item(49, 4).Value = "Montréal"
$ws.Cells.Item(49, 5).Value = "Montréal"
$ws.Cells.Item(49, 6).Value = "A"
$ws.Cells.Item(49, 7).Value = "Cote 1"
$ws.Cells.Item(49, 8).Value = "Donnée soumise à l'AARQ"
$ws.Cells.Item(49, 9).Value = "Pierre-Alexandre Bourgeois"

# Row 51 (new sighting, 2023-05-07 / Montcalm)
$ws.Cells.Item(51, 1).Value = "5/7/2023"
$ws.Cells.Item(51, 2).Value = "BUAM"
$ws.Cells.Item(51, 3).Value = "N/A"
$ws.Cells.Item(51, 4).Value = "Montcalm"
$ws.Cells.Item(51, 5).Value = "Lanaudière"
$ws.Cells.Item(51, 6).Value = "B"
$ws.Cells.Item(51, 7).Value = "Cote 3"
$ws.Cells.Item(51, 8).Value = "Amplexus et ponte, donnée soumise à l'AARQ"
$ws.Cells.Item(51, 9).Value = "Marie-Eve Nepveu"

# Match the saved selection state from the authored workbook.
$ws.Range("E59").Select()
